$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = "kekest"
$ws.Range("D3").Value = "wuest"
$ws.Range("H3").Value = 6.0
$ws.Range("G8").Value = "location 6"
$ws.Range("H8").Value = "name 6"
